$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" column (C) date serial was bumped by one day (2025-03-13 -> 2025-03-14)
# for every data row (rows 2 through 43).
$ws.Range("C2:C43").Value = 45730
